$wb = $excel.ActiveWorkbook

# ---- Sheet: detail_sp ----
$ws = $wb.Worksheets.Item("detail_sp")
$ws.Range("C2").Value = 'American Indian or Alaskan Native'
$ws.Range("C3").Value = 'Asian'
$ws.Range("C4").Value = 'Black or African American'
$ws.Range("C6").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C7").Value = 'Some Other Race'
$ws.Range("C9").Value = 'White'
$ws.Range("C11").Value = 'Multirace PSRC'
$ws.Range("C12").Value = 'Single race PSRC'
$ws.Range("C13").Value = 'Single race Harvard'
$ws.Range("C14").Value = 'People of color'
$ws.Range("C15").Value = 'American Indian or Alaskan Native'
$ws.Range("C16").Value = 'Asian'
$ws.Range("C17").Value = 'Black or African American'
$ws.Range("C19").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C20").Value = 'Some Other Race'
$ws.Range("C22").Value = 'White'
$ws.Range("C24").Value = 'Multirace PSRC'
$ws.Range("C25").Value = 'Single race PSRC'
$ws.Range("C26").Value = 'Single race Harvard'
$ws.Range("C27").Value = 'People of color'
$ws.Range("C28").Value = 'American Indian or Alaskan Native'
$ws.Range("C29").Value = 'Asian'
$ws.Range("C30").Value = 'Black or African American'
$ws.Range("C32").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C33").Value = 'Some Other Race'
$ws.Range("C35").Value = 'White'
$ws.Range("C37").Value = 'Multirace PSRC'
$ws.Range("C38").Value = 'Single race PSRC'
$ws.Range("C39").Value = 'Single race Harvard'
$ws.Range("C40").Value = 'People of color'
$ws.Range("C41").Value = 'American Indian or Alaskan Native'
$ws.Range("C42").Value = 'Asian'
$ws.Range("C43").Value = 'Black or African American'
$ws.Range("C45").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C46").Value = 'Some Other Race'
$ws.Range("C48").Value = 'White'
$ws.Range("C50").Value = 'Multirace PSRC'
$ws.Range("C51").Value = 'Single race PSRC'
$ws.Range("C52").Value = 'Single race Harvard'
$ws.Range("C53").Value = 'People of color'
$ws.Range("C54").Value = 'American Indian or Alaskan Native'
$ws.Range("C55").Value = 'Asian'
$ws.Range("C56").Value = 'Black or African American'
$ws.Range("C58").Value = 'Some Other Race'
$ws.Range("C60").Value = 'White'
$ws.Range("C62").Value = 'Multirace PSRC'
$ws.Range("C63").Value = 'Single race PSRC'
$ws.Range("C64").Value = 'Single race Harvard'
$ws.Range("C65").Value = 'People of color'

# ---- Sheet: detail_mp ----
$ws = $wb.Worksheets.Item("detail_mp")
$ws.Range("C2").Value = 'American Indian or Alaskan Native'
$ws.Range("C3").Value = 'Asian'
$ws.Range("C4").Value = 'Black or African American'
$ws.Range("C6").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C7").Value = 'Some Other Race'
$ws.Range("C9").Value = 'White'
$ws.Range("C11").Value = 'MNAW'
$ws.Range("F11").Value = 10477
$ws.Range("G11").Value = 11683
$ws.Range("I11").Value = 0.446666098226467
$ws.Range("J11").Value = 0.420403022670025
$ws.Range("L11").Value = 1064.65703893437
$ws.Range("M11").Value = 1094.40101879921
$ws.Range("O11").Value = 0.0346965442905105
$ws.Range("P11").Value = 0.0323866693668837
$ws.Range("C12").Value = 'Multirace incl. Asian'
$ws.Range("F12").Value = 13972
$ws.Range("G12").Value = 15596
$ws.Range("I12").Value = 0.604351399281976
$ws.Range("J12").Value = 0.609171158503242
$ws.Range("L12").Value = 1068.3468932257
$ws.Range("M12").Value = 1142.30966385707
$ws.Range("O12").Value = 0.0315685219912864
$ws.Range("P12").Value = 0.0285285953578565
$ws.Range("C13").Value = 'Multirace incl. Asian, white'
$ws.Range("F13").Value = 44421
$ws.Range("G13").Value = 46176
$ws.Range("I13").Value = 0.722398399765819
$ws.Range("J13").Value = 0.726803393512033
$ws.Range("L13").Value = 1786.35839856291
$ws.Range("M13").Value = 1803.61355079428
$ws.Range("O13").Value = 0.0188298659302052
$ws.Range("P13").Value = 0.0175137776059748
$ws.Range("C14").Value = 'Multirace incl. white'
$ws.Range("F14").Value = 103745
$ws.Range("G14").Value = 114796
$ws.Range("I14").Value = 0.635034798523588
$ws.Range("J14").Value = 0.626927130834639
$ws.Range("L14").Value = 3185.33156944316
$ws.Range("M14").Value = 3339.30291863362
$ws.Range("O14").Value = 0.0127117083843686
$ws.Range("P14").Value = 0.0123166870756028
$ws.Range("C15").Value = 'Multirace PSRC'
$ws.Range("C16").Value = 'Single race PSRC'
$ws.Range("C17").Value = 'Multirace Harvard'
$ws.Range("C18").Value = 'Single race Harvard'
$ws.Range("C19").Value = 'People of color'
$ws.Range("C20").Value = 'American Indian or Alaskan Native'
$ws.Range("C21").Value = 'Asian'
$ws.Range("C22").Value = 'Black or African American'
$ws.Range("C24").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C25").Value = 'Some Other Race'
$ws.Range("C27").Value = 'White'
$ws.Range("C29").Value = 'MNAW'
$ws.Range("F29").Value = 4656
$ws.Range("G29").Value = 5120
$ws.Range("I29").Value = 0.402003108271456
$ws.Range("J29").Value = 0.37347727770078
$ws.Range("L29").Value = 632.787583347485
$ws.Range("M29").Value = 670.335077909735
$ws.Range("O29").Value = 0.0543405419322603
$ws.Range("P29").Value = 0.051021618622519
$ws.Range("C30").Value = 'Multirace incl. Asian'
$ws.Range("F30").Value = 7666
$ws.Range("G30").Value = 8749
$ws.Range("I30").Value = 0.566467154363408
$ws.Range("J30").Value = 0.581097236981934
$ws.Range("L30").Value = 858.749794613513
$ws.Range("M30").Value = 928.327442655607
$ws.Range("O30").Value = 0.0438673202832475
$ws.Range("P30").Value = 0.041504622251493
$ws.Range("C31").Value = 'Multirace incl. Asian, white'
$ws.Range("F31").Value = 26737
$ws.Range("G31").Value = 27538
$ws.Range("I31").Value = 0.68079851297329
$ws.Range("J31").Value = 0.684275916906868
$ws.Range("L31").Value = 1318.94295575709
$ws.Range("M31").Value = 1357.66953530029
$ws.Range("O31").Value = 0.0246179144730839
$ws.Range("P31").Value = 0.023457571977104
$ws.Range("C32").Value = 'Multirace incl. white'
$ws.Range("F32").Value = 45254
$ws.Range("G32").Value = 49467
$ws.Range("I32").Value = 0.586184116786052
$ws.Range("J32").Value = 0.584474508182194
$ws.Range("L32").Value = 1994.28425512827
$ws.Range("M32").Value = 2100.50126652426
$ws.Range("O32").Value = 0.0188157743667248
$ws.Range("P32").Value = 0.0176542334654063
$ws.Range("C33").Value = 'Multirace PSRC'
$ws.Range("C34").Value = 'Single race PSRC'
$ws.Range("C35").Value = 'Multirace Harvard'
$ws.Range("C36").Value = 'Single race Harvard'
$ws.Range("C37").Value = 'People of color'
$ws.Range("C38").Value = 'American Indian or Alaskan Native'
$ws.Range("C39").Value = 'Asian'
$ws.Range("C40").Value = 'Black or African American'
$ws.Range("C42").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C43").Value = 'Some Other Race'
$ws.Range("C45").Value = 'White'
$ws.Range("C47").Value = 'MNAW'
$ws.Range("F47").Value = 648
$ws.Range("G47").Value = 657
$ws.Range("I47").Value = 0.374133949191686
$ws.Range("J47").Value = 0.352278820375335
$ws.Range("L47").Value = 197.144516963445
$ws.Range("M47").Value = 196.631178863501
$ws.Range("O47").Value = 0.102877942826686
$ws.Range("P47").Value = 0.099943556865988
$ws.Range("C48").Value = 'Multirace incl. Asian'
$ws.Range("F48").Value = 504
$ws.Range("G48").Value = 618
$ws.Range("I48").Value = 0.609431680773882
$ws.Range("J48").Value = 0.647120418848168
$ws.Range("L48").Value = 161.698563251193
$ws.Range("M48").Value = 181.973733633456
$ws.Range("O48").Value = 0.141283767149141
$ws.Range("P48").Value = 0.129834215880133
$ws.Range("R48").Value = 'fair'
$ws.Range("S48").Value = 'fair'
$ws.Range("C49").Value = 'Multirace incl. Asian, white'
$ws.Range("F49").Value = 2968
$ws.Range("G49").Value = 3027
$ws.Range("I49").Value = 0.834880450070324
$ws.Range("J49").Value = 0.837576092971776
$ws.Range("L49").Value = 523.929271968078
$ws.Range("M49").Value = 525.725503417658
$ws.Range("O49").Value = 0.0600551876947096
$ws.Range("P49").Value = 0.0591642615070534
$ws.Range("C50").Value = 'Multirace incl. white'
$ws.Range("F50").Value = 8515
$ws.Range("G50").Value = 9595
$ws.Range("I50").Value = 0.681363527246539
$ws.Range("J50").Value = 0.668222021032105
$ws.Range("L50").Value = 735.432201509765
$ws.Range("M50").Value = 776.316353767103
$ws.Range("O50").Value = 0.040151419122855
$ws.Range("P50").Value = 0.0377222647256608
$ws.Range("R50").Value = 'good'
$ws.Range("S50").Value = 'good'
$ws.Range("C51").Value = 'Multirace PSRC'
$ws.Range("C52").Value = 'Single race PSRC'
$ws.Range("C53").Value = 'Multirace Harvard'
$ws.Range("C54").Value = 'Single race Harvard'
$ws.Range("C55").Value = 'People of color'
$ws.Range("C56").Value = 'American Indian or Alaskan Native'
$ws.Range("C57").Value = 'Asian'
$ws.Range("C58").Value = 'Black or African American'
$ws.Range("C60").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C61").Value = 'Some Other Race'
$ws.Range("C63").Value = 'White'
$ws.Range("C65").Value = 'MNAW'
$ws.Range("F65").Value = 3816
$ws.Range("G65").Value = 4334
$ws.Range("I65").Value = 0.531994981179423
$ws.Range("J65").Value = 0.496505899873983
$ws.Range("L65").Value = 730.175643145196
$ws.Range("M65").Value = 733.935604669272
$ws.Range("O65").Value = 0.07019852353216
$ws.Range("P65").Value = 0.0597673609065604
$ws.Range("C66").Value = 'Multirace incl. Asian'
$ws.Range("F66").Value = 3370
$ws.Range("G66").Value = 3615
$ws.Range("I66").Value = 0.64509954058193
$ws.Range("J66").Value = 0.622096024780589
$ws.Range("L66").Value = 486.088175460482
$ws.Range("M66").Value = 497.215230262509
$ws.Range("O66").Value = 0.0712422612573767
$ws.Range("P66").Value = 0.0654447160060862
$ws.Range("C67").Value = 'Multirace incl. Asian, white'
$ws.Range("F67").Value = 7363
$ws.Range("G67").Value = 7728
$ws.Range("I67").Value = 0.758759274525969
$ws.Range("J67").Value = 0.764997030291031
$ws.Range("L67").Value = 722.108346442901
$ws.Range("M67").Value = 770.140442510326
$ws.Range("O67").Value = 0.0524766273781432
$ws.Range("P67").Value = 0.0512363326372591
$ws.Range("C68").Value = 'Multirace incl. white'
$ws.Range("F68").Value = 28809
$ws.Range("G68").Value = 32197
$ws.Range("I68").Value = 0.669883272101567
$ws.Range("J68").Value = 0.658371503353509
$ws.Range("L68").Value = 1289.23368504948
$ws.Range("M68").Value = 1393.86121937946
$ws.Range("O68").Value = 0.0219340958081491
$ws.Range("P68").Value = 0.0215379225275967
$ws.Range("C69").Value = 'Multirace PSRC'
$ws.Range("C70").Value = 'Single race PSRC'
$ws.Range("C71").Value = 'Multirace Harvard'
$ws.Range("C72").Value = 'Single race Harvard'
$ws.Range("C73").Value = 'People of color'
$ws.Range("C74").Value = 'American Indian or Alaskan Native'
$ws.Range("C75").Value = 'Asian'
$ws.Range("C76").Value = 'Black or African American'
$ws.Range("C78").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C79").Value = 'Some Other Race'
$ws.Range("C81").Value = 'White'
$ws.Range("C83").Value = 'MNAW'
$ws.Range("F83").Value = 1357
$ws.Range("G83").Value = 1572
$ws.Range("I83").Value = 0.457056247894914
$ws.Range("J83").Value = 0.450817321479782
$ws.Range("L83").Value = 311.283815516483
$ws.Range("M83").Value = 352.713472027934
$ws.Range("O83").Value = 0.0953887538322556
$ws.Range("P83").Value = 0.0883245661055159
$ws.Range("C84").Value = 'Multirace incl. Asian'
$ws.Range("F84").Value = 2432
$ws.Range("G84").Value = 2614
$ws.Range("I84").Value = 0.687977369165488
$ws.Range("J84").Value = 0.691534391534392
$ws.Range("L84").Value = 391.677747151597
$ws.Range("M84").Value = 402.444635703161
$ws.Range("O84").Value = 0.0767110124360402
$ws.Range("P84").Value = 0.0743627692865157
$ws.Range("C85").Value = 'Multirace incl. Asian, white'
$ws.Range("F85").Value = 7353
$ws.Range("G85").Value = 7883
$ws.Range("I85").Value = 0.82073892175466
$ws.Range("J85").Value = 0.823461819701243
$ws.Range("L85").Value = 750.004552519516
$ws.Range("M85").Value = 792.229164965858
$ws.Range("O85").Value = 0.038345665678746
$ws.Range("P85").Value = 0.03738176549344
$ws.Range("C86").Value = 'Multirace incl. white'
$ws.Range("F86").Value = 21167
$ws.Range("G86").Value = 23537
$ws.Range("I86").Value = 0.690265775313876
$ws.Range("J86").Value = 0.66845588026469
$ws.Range("L86").Value = 1163.54746721825
$ws.Range("M86").Value = 1279.93754457639
$ws.Range("O86").Value = 0.0279240792112708
$ws.Range("P86").Value = 0.0244737045498222
$ws.Range("C87").Value = 'Multirace PSRC'
$ws.Range("C88").Value = 'Single race PSRC'
$ws.Range("C89").Value = 'Multirace Harvard'
$ws.Range("C90").Value = 'Single race Harvard'
$ws.Range("C91").Value = 'People of color'

# ---- Sheet: dichot_sp ----
$ws = $wb.Worksheets.Item("dichot_sp")
$ws.Range("C2").Value = 'American Indian or Alaskan Native'
$ws.Range("C3").Value = 'Asian'
$ws.Range("C4").Value = 'Black or African American'
$ws.Range("C6").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C7").Value = 'Some Other Race'
$ws.Range("C9").Value = 'White'
$ws.Range("C11").Value = 'Multirace PSRC'
$ws.Range("C12").Value = 'Single race PSRC'
$ws.Range("C13").Value = 'Single race Harvard'
$ws.Range("C14").Value = 'People of color'
$ws.Range("C15").Value = 'American Indian or Alaskan Native'
$ws.Range("C16").Value = 'Asian'
$ws.Range("C17").Value = 'Black or African American'
$ws.Range("C19").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C20").Value = 'Some Other Race'
$ws.Range("C22").Value = 'White'
$ws.Range("C24").Value = 'Multirace PSRC'
$ws.Range("C25").Value = 'Single race PSRC'
$ws.Range("C26").Value = 'Single race Harvard'
$ws.Range("C27").Value = 'People of color'
$ws.Range("C28").Value = 'American Indian or Alaskan Native'
$ws.Range("C29").Value = 'Asian'
$ws.Range("C30").Value = 'Black or African American'
$ws.Range("C32").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C33").Value = 'Some Other Race'
$ws.Range("C35").Value = 'White'
$ws.Range("C37").Value = 'Multirace PSRC'
$ws.Range("C38").Value = 'Single race PSRC'
$ws.Range("C39").Value = 'Single race Harvard'
$ws.Range("C40").Value = 'People of color'
$ws.Range("C41").Value = 'American Indian or Alaskan Native'
$ws.Range("C42").Value = 'Asian'
$ws.Range("C43").Value = 'Black or African American'
$ws.Range("C45").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C46").Value = 'Some Other Race'
$ws.Range("C48").Value = 'White'
$ws.Range("C50").Value = 'Multirace PSRC'
$ws.Range("C51").Value = 'Single race PSRC'
$ws.Range("C52").Value = 'Single race Harvard'
$ws.Range("C53").Value = 'People of color'
$ws.Range("C54").Value = 'American Indian or Alaskan Native'
$ws.Range("C55").Value = 'Asian'
$ws.Range("C56").Value = 'Black or African American'
$ws.Range("C58").Value = 'Some Other Race'
$ws.Range("C60").Value = 'White'
$ws.Range("C62").Value = 'Multirace PSRC'
$ws.Range("C63").Value = 'Single race PSRC'
$ws.Range("C64").Value = 'Single race Harvard'
$ws.Range("C65").Value = 'People of color'

# ---- Sheet: dichot_mp ----
$ws = $wb.Worksheets.Item("dichot_mp")
$ws.Range("C2").Value = 'American Indian or Alaskan Native'
$ws.Range("C3").Value = 'Asian'
$ws.Range("C4").Value = 'Black or African American'
$ws.Range("C6").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C7").Value = 'Some Other Race'
$ws.Range("C9").Value = 'White'
$ws.Range("C11").Value = 'MNW'
$ws.Range("F11").Value = 24449
$ws.Range("G11").Value = 27279
$ws.Range("I11").Value = 0.524938271604938
$ws.Range("J11").Value = 0.510919238837279
$ws.Range("L11").Value = 1388.41712894666
$ws.Range("M11").Value = 1612.84964356105
$ws.Range("O11").Value = 0.023210589240122
$ws.Range("P11").Value = 0.0227267115378057
$ws.Range("C12").Value = 'Multirace incl. white'
$ws.Range("F12").Value = 148166
$ws.Range("G12").Value = 160972
$ws.Range("I12").Value = 0.658925553677844
$ws.Range("J12").Value = 0.652654454634653
$ws.Range("L12").Value = 3823.40580761185
$ws.Range("M12").Value = 4052.37436437964
$ws.Range("O12").Value = 0.0114449632085157
$ws.Range("P12").Value = 0.0109241531840316
$ws.Range("C13").Value = 'Multirace PSRC'
$ws.Range("C14").Value = 'Single race PSRC'
$ws.Range("C15").Value = 'Multirace Harvard'
$ws.Range("C16").Value = 'Single race Harvard'
$ws.Range("C17").Value = 'People of color'
$ws.Range("C18").Value = 'American Indian or Alaskan Native'
$ws.Range("C19").Value = 'Asian'
$ws.Range("C20").Value = 'Black or African American'
$ws.Range("C22").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C23").Value = 'Some Other Race'
$ws.Range("C25").Value = 'White'
$ws.Range("C27").Value = 'MNW'
$ws.Range("F27").Value = 12322
$ws.Range("G27").Value = 13869
$ws.Range("I27").Value = 0.490623133585507
$ws.Range("J27").Value = 0.482148444289936
$ws.Range("L27").Value = 953.507048505673
$ws.Range("M27").Value = 1075.61171351992
$ws.Range("O27").Value = 0.032759640243789
$ws.Range("P27").Value = 0.0340526388540633
$ws.Range("C28").Value = 'Multirace incl. white'
$ws.Range("F28").Value = 71991
$ws.Range("G28").Value = 77005
$ws.Range("I28").Value = 0.618086439892165
$ws.Range("J28").Value = 0.616636904523579
$ws.Range("L28").Value = 2437.93659123833
$ws.Range("M28").Value = 2604.7269407603
$ws.Range("O28").Value = 0.0155289782052378
$ws.Range("P28").Value = 0.0143730173722226
$ws.Range("C29").Value = 'Multirace PSRC'
$ws.Range("C30").Value = 'Single race PSRC'
$ws.Range("C31").Value = 'Multirace Harvard'
$ws.Range("C32").Value = 'Single race Harvard'
$ws.Range("C33").Value = 'People of color'
$ws.Range("C34").Value = 'American Indian or Alaskan Native'
$ws.Range("C35").Value = 'Asian'
$ws.Range("C36").Value = 'Black or African American'
$ws.Range("C38").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C39").Value = 'Some Other Race'
$ws.Range("C41").Value = 'White'
$ws.Range("C43").Value = 'MNW'
$ws.Range("F43").Value = 1152
$ws.Range("G43").Value = 1275
$ws.Range("I43").Value = 0.450175849941383
$ws.Range("J43").Value = 0.452127659574468
$ws.Range("L43").Value = 247.323165142794
$ws.Range("M43").Value = 258.711749750084
$ws.Range("O43").Value = 0.0868024062370842
$ws.Range("P43").Value = 0.0864185452800276
$ws.Range("C44").Value = 'Multirace incl. white'
$ws.Range("F44").Value = 11483
$ws.Range("G44").Value = 12622
$ws.Range("I44").Value = 0.715362571642163
$ws.Range("J44").Value = 0.702275635675736
$ws.Range("L44").Value = 864.874296400928
$ws.Range("M44").Value = 851.158194769045
$ws.Range("O44").Value = 0.0332561989326599
$ws.Range("P44").Value = 0.0311987614724177
$ws.Range("C45").Value = 'Multirace PSRC'
$ws.Range("C46").Value = 'Single race PSRC'
$ws.Range("C47").Value = 'Multirace Harvard'
$ws.Range("C48").Value = 'Single race Harvard'
$ws.Range("C49").Value = 'People of color'
$ws.Range("C50").Value = 'American Indian or Alaskan Native'
$ws.Range("C51").Value = 'Asian'
$ws.Range("C52").Value = 'Black or African American'
$ws.Range("C54").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C55").Value = 'Some Other Race'
$ws.Range("C57").Value = 'White'
$ws.Range("C59").Value = 'MNW'
$ws.Range("F59").Value = 7186
$ws.Range("G59").Value = 7949
$ws.Range("I59").Value = 0.579656368476244
$ws.Range("J59").Value = 0.546698762035763
$ws.Range("L59").Value = 835.296920139779
$ws.Range("M59").Value = 848.412398016878
$ws.Range("O59").Value = 0.0490734738141149
$ws.Range("P59").Value = 0.0441927448472436
$ws.Range("C60").Value = 'Multirace incl. white'
$ws.Range("F60").Value = 36172
$ws.Range("G60").Value = 39925
$ws.Range("I60").Value = 0.686245494213622
$ws.Range("J60").Value = 0.676626105819747
$ws.Range("L60").Value = 1508.09761079522
$ws.Range("M60").Value = 1595.23038503221
$ws.Range("O60").Value = 0.0214697141296225
$ws.Range("P60").Value = 0.0204485645550855
$ws.Range("C61").Value = 'Multirace PSRC'
$ws.Range("C62").Value = 'Single race PSRC'
$ws.Range("C63").Value = 'Multirace Harvard'
$ws.Range("C64").Value = 'Single race Harvard'
$ws.Range("C65").Value = 'People of color'
$ws.Range("C66").Value = 'American Indian or Alaskan Native'
$ws.Range("C67").Value = 'Asian'
$ws.Range("C68").Value = 'Black or African American'
$ws.Range("C70").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C71").Value = 'Some Other Race'
$ws.Range("C73").Value = 'White'
$ws.Range("C75").Value = 'MNW'
$ws.Range("F75").Value = 3789
$ws.Range("G75").Value = 4186
$ws.Range("I75").Value = 0.582564575645756
$ws.Range("J75").Value = 0.57602862254025
$ws.Range("L75").Value = 505.071393024046
$ws.Range("M75").Value = 561.655975752729
$ws.Range("O75").Value = 0.0636070345002761
$ws.Range("P75").Value = 0.0620669383785101
$ws.Range("C76").Value = 'Multirace incl. white'
$ws.Range("F76").Value = 28520
$ws.Range("G76").Value = 31420
$ws.Range("I76").Value = 0.719765798505956
$ws.Range("J76").Value = 0.701589853519114
$ws.Range("L76").Value = 1491.74020347043
$ws.Range("M76").Value = 1686.19645568102
$ws.Range("O76").Value = 0.0246477713919709
$ws.Range("P76").Value = 0.0219708181945032
$ws.Range("C77").Value = 'Multirace PSRC'
$ws.Range("C78").Value = 'Single race PSRC'
$ws.Range("C79").Value = 'Multirace Harvard'
$ws.Range("C80").Value = 'Single race Harvard'
$ws.Range("C81").Value = 'People of color'

# ---- Sheet: single_sp ----
$ws = $wb.Worksheets.Item("single_sp")
$ws.Range("C2").Value = 'American Indian or Alaskan Native'
$ws.Range("C3").Value = 'Asian'
$ws.Range("C4").Value = 'Black or African American'
$ws.Range("C6").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C7").Value = 'Some Other Race'
$ws.Range("C9").Value = 'White'
$ws.Range("C11").Value = 'Multirace PSRC'
$ws.Range("C12").Value = 'Single race PSRC'
$ws.Range("C13").Value = 'Single race Harvard'
$ws.Range("C14").Value = 'People of color'
$ws.Range("C15").Value = 'American Indian or Alaskan Native'
$ws.Range("C16").Value = 'Asian'
$ws.Range("C17").Value = 'Black or African American'
$ws.Range("C19").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C20").Value = 'Some Other Race'
$ws.Range("C22").Value = 'White'
$ws.Range("C24").Value = 'Multirace PSRC'
$ws.Range("C25").Value = 'Single race PSRC'
$ws.Range("C26").Value = 'Single race Harvard'
$ws.Range("C27").Value = 'People of color'
$ws.Range("C28").Value = 'American Indian or Alaskan Native'
$ws.Range("C29").Value = 'Asian'
$ws.Range("C30").Value = 'Black or African American'
$ws.Range("C32").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C33").Value = 'Some Other Race'
$ws.Range("C35").Value = 'White'
$ws.Range("C37").Value = 'Multirace PSRC'
$ws.Range("C38").Value = 'Single race PSRC'
$ws.Range("C39").Value = 'Single race Harvard'
$ws.Range("C40").Value = 'People of color'
$ws.Range("C41").Value = 'American Indian or Alaskan Native'
$ws.Range("C42").Value = 'Asian'
$ws.Range("C43").Value = 'Black or African American'
$ws.Range("C45").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C46").Value = 'Some Other Race'
$ws.Range("C48").Value = 'White'
$ws.Range("C50").Value = 'Multirace PSRC'
$ws.Range("C51").Value = 'Single race PSRC'
$ws.Range("C52").Value = 'Single race Harvard'
$ws.Range("C53").Value = 'People of color'
$ws.Range("C54").Value = 'American Indian or Alaskan Native'
$ws.Range("C55").Value = 'Asian'
$ws.Range("C56").Value = 'Black or African American'
$ws.Range("C58").Value = 'Some Other Race'
$ws.Range("C60").Value = 'White'
$ws.Range("C62").Value = 'Multirace PSRC'
$ws.Range("C63").Value = 'Single race PSRC'
$ws.Range("C64").Value = 'Single race Harvard'
$ws.Range("C65").Value = 'People of color'

# ---- Sheet: single_mp ----
$ws = $wb.Worksheets.Item("single_mp")
$ws.Range("C2").Value = 'American Indian or Alaskan Native'
$ws.Range("C3").Value = 'Asian'
$ws.Range("C4").Value = 'Black or African American'
$ws.Range("C6").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C7").Value = 'Some Other Race'
$ws.Range("C9").Value = 'White'
$ws.Range("C11").Value = 'Multirace'
$ws.Range("C12").Value = 'Multirace PSRC'
$ws.Range("C13").Value = 'Single race PSRC'
$ws.Range("C14").Value = 'Multirace Harvard'
$ws.Range("C15").Value = 'Single race Harvard'
$ws.Range("C16").Value = 'People of color'
$ws.Range("C17").Value = 'American Indian or Alaskan Native'
$ws.Range("C18").Value = 'Asian'
$ws.Range("C19").Value = 'Black or African American'
$ws.Range("C21").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C22").Value = 'Some Other Race'
$ws.Range("C24").Value = 'White'
$ws.Range("C26").Value = 'Multirace'
$ws.Range("C27").Value = 'Multirace PSRC'
$ws.Range("C28").Value = 'Single race PSRC'
$ws.Range("C29").Value = 'Multirace Harvard'
$ws.Range("C30").Value = 'Single race Harvard'
$ws.Range("C31").Value = 'People of color'
$ws.Range("C32").Value = 'American Indian or Alaskan Native'
$ws.Range("C33").Value = 'Asian'
$ws.Range("C34").Value = 'Black or African American'
$ws.Range("C36").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C37").Value = 'Some Other Race'
$ws.Range("C39").Value = 'White'
$ws.Range("C41").Value = 'Multirace'
$ws.Range("C42").Value = 'Multirace PSRC'
$ws.Range("C43").Value = 'Single race PSRC'
$ws.Range("C44").Value = 'Multirace Harvard'
$ws.Range("C45").Value = 'Single race Harvard'
$ws.Range("C46").Value = 'People of color'
$ws.Range("C47").Value = 'American Indian or Alaskan Native'
$ws.Range("C48").Value = 'Asian'
$ws.Range("C49").Value = 'Black or African American'
$ws.Range("C51").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C52").Value = 'Some Other Race'
$ws.Range("C54").Value = 'White'
$ws.Range("C56").Value = 'Multirace'
$ws.Range("C57").Value = 'Multirace PSRC'
$ws.Range("C58").Value = 'Single race PSRC'
$ws.Range("C59").Value = 'Multirace Harvard'
$ws.Range("C60").Value = 'Single race Harvard'
$ws.Range("C61").Value = 'People of color'
$ws.Range("C62").Value = 'American Indian or Alaskan Native'
$ws.Range("C63").Value = 'Asian'
$ws.Range("C64").Value = 'Black or African American'
$ws.Range("C66").Value = 'Native Hawaiian or Pacific Islander'
$ws.Range("C67").Value = 'Some Other Race'
$ws.Range("C69").Value = 'White'
$ws.Range("C71").Value = 'Multirace'
$ws.Range("C72").Value = 'Multirace PSRC'
$ws.Range("C73").Value = 'Single race PSRC'
$ws.Range("C74").Value = 'Multirace Harvard'
$ws.Range("C75").Value = 'Single race Harvard'
$ws.Range("C76").Value = 'People of color'
